$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3030883333333333
$ws.Range("H2").Value = 0.909265
$ws.Range("I2").Value = 0.5850568929085261
$ws.Range("J2").Value = 0.585056892908526
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 21.19112166666667
$ws.Range("N2").Value = 63.573365
$ws.Range("O2").Value = 0.4896103362399876
$ws.Range("P2").Value = 0.4896103362399877
$ws.Range("Q2").Value = 6.42278174741389
$ws.Range("R2").Value = 57.805035726725
$ws.Range("S2").Value = 0.2864499020564659
$ws.Range("T2").Value = 0.2864499020564659
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3030883333333333
$ws.Range("H3").Value = 0.909265
$ws.Range("I3").Value = 0.5850568929085261
$ws.Range("J3").Value = 0.585056892908526
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 13.462409
$ws.Range("N3").Value = 40.387227
$ws.Range("O3").Value = 0.3110422704739745
$ws.Range("P3").Value = 0.3110422704739746
$ws.Range("Q3").Value = 4.080299106461666
$ws.Range("R3").Value = 36.722691958155
$ws.Range("S3").Value = 0.1819774243267169
$ws.Range("T3").Value = 0.1819774243267169
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3030883333333333
$ws.Range("H4").Value = 0.909265
$ws.Range("I4").Value = 0.5850568929085261
$ws.Range("J4").Value = 0.585056892908526
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.601420333333333
$ws.Range("N4").Value = 25.804261
$ws.Range("O4").Value = 0.1987315427559073
$ws.Range("P4").Value = 0.1987315427559073
$ws.Range("Q4").Value = 2.606990153129445
$ws.Range("R4").Value = 23.462911378165
$ws.Range("S4").Value = 0.116269258927689
$ws.Range("T4").Value = 0.116269258927689
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3030883333333333
$ws.Range("H5").Value = 0.909265
$ws.Range("I5").Value = 0.5850568929085261
$ws.Range("J5").Value = 0.585056892908526
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.026655
$ws.Range("N5").Value = 0.079965
$ws.Range("O5").Value = 0.0006158505301305131
$ws.Range("P5").Value = 0.0006158505301305133
$ws.Range("Q5").Value = 0.008078819525
$ws.Range("R5").Value = 0.072709375725
$ws.Range("S5").Value = 0.0003603075976542267
$ws.Range("T5").Value = 0.0003603075976542267
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.214961
$ws.Range("H6").Value = 0.644883
$ws.Range("I6").Value = 0.4149431070914739
$ws.Range("J6").Value = 0.4149431070914739
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 21.19112166666667
$ws.Range("N6").Value = 63.573365
$ws.Range("O6").Value = 0.4896103362399876
$ws.Range("P6").Value = 0.4896103362399877
$ws.Range("Q6").Value = 4.555264704588333
$ws.Range("R6").Value = 40.997382341295
$ws.Range("S6").Value = 0.2031604341835217
$ws.Range("T6").Value = 0.2031604341835218
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.214961
$ws.Range("H7").Value = 0.644883
$ws.Range("I7").Value = 0.4149431070914739
$ws.Range("J7").Value = 0.4149431070914739
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.462409
$ws.Range("N7").Value = 40.387227
$ws.Range("O7").Value = 0.3110422704739745
$ws.Range("P7").Value = 0.3110422704739746
$ws.Range("Q7").Value = 2.893892901049
$ws.Range("R7").Value = 26.045036109441
$ws.Range("S7").Value = 0.1290648461472576
$ws.Range("T7").Value = 0.1290648461472576
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.214961
$ws.Range("H8").Value = 0.644883
$ws.Range("I8").Value = 0.4149431070914739
$ws.Range("J8").Value = 0.4149431070914739
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.601420333333333
$ws.Range("N8").Value = 25.804261
$ws.Range("O8").Value = 0.1987315427559073
$ws.Range("P8").Value = 0.1987315427559073
$ws.Range("Q8").Value = 1.848969916273667
$ws.Range("R8").Value = 16.640729246463
$ws.Range("S8").Value = 0.08246228382821824
$ws.Range("T8").Value = 0.08246228382821827
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.214961
$ws.Range("H9").Value = 0.644883
$ws.Range("I9").Value = 0.4149431070914739
$ws.Range("J9").Value = 0.4149431070914739
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.026655
$ws.Range("N9").Value = 0.079965
$ws.Range("O9").Value = 0.0006158505301305131
$ws.Range("P9").Value = 0.0006158505301305133
$ws.Range("Q9").Value = 0.005729785454999999
$ws.Range("R9").Value = 0.051568069095
$ws.Range("S9").Value = 0.0002555429324762864
$ws.Range("T9").Value = 0.0002555429324762865
